$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — set the text first, then copy the
# existing header formatting (bold/border/centered, style index used by the
# rest of row 1) over via a formats-only paste so no new style entries are
# minted.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns I ("I0") and J ("IF") — same value per row.
$values = @{
    2  = 8
    3  = 8
    4  = 9
    5  = 7
    6  = 1
    7  = 2
    8  = 4
    9  = 9
    10 = 9
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v   # column I
    $ws.Cells.Item($row, 10).Value = $v  # column J
}
